$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure percent-like text values (e.g. '62%') are set as NumberFormat text
# so Excel doesn't auto-convert them to numeric percentages.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H45").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-25 02:18:19"
$ws.Range("H2").Value = "62%"
$ws.Range("O2").Value = "1.9 °C"
$ws.Range("E3").Value = "2026-02-25 02:18:21"
$ws.Range("H3").Value = "35%"
$ws.Range("O3").Value = "2.2 °C"
$ws.Range("E4").Value = "2026-02-25 02:18:23"
$ws.Range("N4").Value = "2.7 °C 1:48 TU"
$ws.Range("O4").Value = "3.3 °C"
$ws.Range("E5").Value = "2026-02-25 02:18:25"
$ws.Range("H5").Value = "24%"
$ws.Range("O5").Value = "5.3 °C"
$ws.Range("E6").Value = "2026-02-25 02:18:27"
$ws.Range("H6").Value = "98%"
$ws.Range("L6").Value = "10.8 km/h - 296º 1:35 TU"
$ws.Range("O6").Value = "9.1 °C"
$ws.Range("E7").Value = "2026-02-25 02:18:30"
$ws.Range("O7").Value = "11.3 °C"
$ws.Range("E8").Value = "2026-02-25 02:18:32"
$ws.Range("M8").Value = "15.5 °C 1:56 TU"
$ws.Range("O8").Value = "15.2 °C"
$ws.Range("E9").Value = "2026-02-25 02:18:34"
$ws.Range("E10").Value = "2026-02-25 02:18:36"
$ws.Range("N10").Value = "3.7 °C 1:59 TU"
$ws.Range("O10").Value = "5.2 °C"
$ws.Range("E11").Value = "2026-02-25 02:18:39"
$ws.Range("H11").Value = "90%"
$ws.Range("N11").Value = "3.3 °C 1:40 TU"
$ws.Range("O11").Value = "3.7 °C"
$ws.Range("E12").Value = "2026-02-25 02:18:41"
$ws.Range("N12").Value = "4.7 °C 1:56 TU"
$ws.Range("O12").Value = "5.3 °C"
$ws.Range("E13").Value = "2026-02-25 02:18:43"
$ws.Range("H13").Value = "90%"
$ws.Range("J13").Value = "1026.2 hPa"
$ws.Range("N13").Value = "-1.7 °C 1:48 TU"
$ws.Range("O13").Value = "-0.5 °C"
$ws.Range("E14").Value = "2026-02-25 02:18:45"
$ws.Range("M14").Value = "7.9 °C 1:50 TU"
$ws.Range("O14").Value = "6.4 °C"
$ws.Range("E15").Value = "2026-02-25 02:18:47"
$ws.Range("N15").Value = "4.4 °C 1:56 TU"
$ws.Range("O15").Value = "5.1 °C"
$ws.Range("E16").Value = "2026-02-25 02:18:49"
$ws.Range("M16").Value = "4.7 °C 1:31 TU"
$ws.Range("E17").Value = "2026-02-25 02:18:52"
$ws.Range("E18").Value = "2026-02-25 02:18:54"
$ws.Range("N18").Value = "6.6 °C 1:57 TU"
$ws.Range("O18").Value = "7.3 °C"
$ws.Range("E19").Value = "2026-02-25 02:18:56"
$ws.Range("H19").Value = "56%"
$ws.Range("E20").Value = "2026-02-25 02:18:58"
$ws.Range("E21").Value = "2026-02-25 02:19:00"
$ws.Range("H21").Value = "72%"
$ws.Range("J21").Value = "1023.1 hPa"
$ws.Range("N21").Value = "3.7 °C 1:58 TU"
$ws.Range("O21").Value = "4.7 °C"
$ws.Range("E22").Value = "2026-02-25 02:19:03"
$ws.Range("O22").Value = "1.0 °C"
$ws.Range("E23").Value = "2026-02-25 02:19:05"
$ws.Range("N23").Value = "2.0 °C 1:56 TU"
$ws.Range("O23").Value = "3.1 °C"
$ws.Range("E24").Value = "2026-02-25 02:19:07"
$ws.Range("N24").Value = "3.6 °C 1:47 TU"
$ws.Range("O24").Value = "4.6 °C"
$ws.Range("E25").Value = "2026-02-25 02:19:09"
$ws.Range("H25").Value = "28%"
$ws.Range("M25").Value = "3.9 °C 1:54 TU"
$ws.Range("O25").Value = "3.5 °C"
$ws.Range("E26").Value = "2026-02-25 02:19:12"
$ws.Range("J26").Value = "1018.9 hPa"
$ws.Range("N26").Value = "8.7 °C 1:45 TU"
$ws.Range("O26").Value = "9.6 °C"
$ws.Range("E27").Value = "2026-02-25 02:19:14"
$ws.Range("G27").Value = "158 cm"
$ws.Range("H27").Value = "41%"
$ws.Range("L27").Value = "21.2 km/h - 239º 1:39 TU"
$ws.Range("M27").Value = "4.3 °C 1:54 TU"
$ws.Range("O27").Value = "3.6 °C"
$ws.Range("E28").Value = "2026-02-25 02:19:16"
$ws.Range("N28").Value = "3.5 °C 1:46 TU"
$ws.Range("O28").Value = "4.5 °C"
$ws.Range("E29").Value = "2026-02-25 02:19:18"
$ws.Range("E30").Value = "2026-02-25 02:19:21"
$ws.Range("H30").Value = "98%"
$ws.Range("N30").Value = "7.8 °C 1:56 TU"
$ws.Range("O30").Value = "8.2 °C"
$ws.Range("E31").Value = "2026-02-25 02:19:23"
$ws.Range("H31").Value = "85%"
$ws.Range("L31").Value = "43.6 km/h - 353º 1:56 TU"
$ws.Range("O31").Value = "10.7 °C"
$ws.Range("E32").Value = "2026-02-25 02:19:25"
$ws.Range("O32").Value = "2.2 °C"
$ws.Range("E33").Value = "2026-02-25 02:19:27"
$ws.Range("H33").Value = "62%"
$ws.Range("J33").Value = "1023.6 hPa"
$ws.Range("N33").Value = "2.5 °C 1:55 TU"
$ws.Range("O33").Value = "3.3 °C"
$ws.Range("E34").Value = "2026-02-25 02:19:29"
$ws.Range("G34").Value = "79 cm"
$ws.Range("H34").Value = "59%"
$ws.Range("N34").Value = "0.6 °C 1:51 TU"
$ws.Range("O34").Value = "2.2 °C"
$ws.Range("E35").Value = "2026-02-25 02:19:32"
$ws.Range("H35").Value = "29%"
$ws.Range("J35").Value = "1019.5 hPa"
$ws.Range("N35").Value = "9.6 °C 1:57 TU"
$ws.Range("O35").Value = "10.8 °C"
$ws.Range("E36").Value = "2026-02-25 02:19:34"
$ws.Range("J36").Value = "1018.8 hPa"
$ws.Range("N36").Value = "7.4 °C 1:46 TU"
$ws.Range("O36").Value = "8.3 °C"
$ws.Range("E37").Value = "2026-02-25 02:19:36"
$ws.Range("H37").Value = "97%"
$ws.Range("J37").Value = "1023.9 hPa"
$ws.Range("N37").Value = "0.9 °C 1:59 TU"
$ws.Range("O37").Value = "1.6 °C"
$ws.Range("E38").Value = "2026-02-25 02:19:39"
$ws.Range("L38").Value = "6.1 km/h - 284º 1:55 TU"
$ws.Range("N38").Value = "3.8 °C 1:50 TU"
$ws.Range("O38").Value = "4.7 °C"
$ws.Range("E39").Value = "2026-02-25 02:19:41"
$ws.Range("H39").Value = "44%"
$ws.Range("M39").Value = "3.0 °C 1:52 TU"
$ws.Range("O39").Value = "2.3 °C"
$ws.Range("E40").Value = "2026-02-25 02:19:43"
$ws.Range("N40").Value = "1.6 °C 1:45 TU"
$ws.Range("O40").Value = "2.2 °C"
$ws.Range("E41").Value = "2026-02-25 02:19:45"
$ws.Range("L41").Value = "4.7 km/h - 30º 1:49 TU"
$ws.Range("E42").Value = "2026-02-25 02:19:47"
$ws.Range("M42").Value = "8.5 °C 1:53 TU"
$ws.Range("O42").Value = "7.9 °C"
$ws.Range("E43").Value = "2026-02-25 02:19:50"
$ws.Range("N43").Value = "3.9 °C 1:59 TU"
$ws.Range("O43").Value = "4.6 °C"
$ws.Range("E44").Value = "2026-02-25 02:19:52"
$ws.Range("H44").Value = "49%"
$ws.Range("O44").Value = "-0.2 °C"
$ws.Range("E45").Value = "2026-02-25 02:19:54"
$ws.Range("H45").Value = "49%"
$ws.Range("N45").Value = "5.1 °C 1:56 TU"
$ws.Range("O45").Value = "6.8 °C"
$ws.Range("E46").Value = "2026-02-25 02:19:57"
$ws.Range("N46").Value = "3.6 °C 1:59 TU"
$ws.Range("O46").Value = "4.5 °C"
